$wb = $excel.ActiveWorkbook

# "想去人数" (want-to-go headcount) column F updates, applied identically
# to both the "展览" and "全部类型" sheets (they mirror the same rows).
$updates = @{
    4  = 1554
    6  = 1084
    7  = 11264
    8  = 11
    10 = 419
    11 = 336
    12 = 1081
    14 = 12290
    15 = 12922
    22 = 74
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
